# Weekly BRVM data refresh (automated) - updates Recommandations + Top_YTD sheets
$wb = $excel.ActiveWorkbook

# --- Sheet "Recommandations": refresh data rows, then drop the now-unused trailing rows ---
$ws1 = $wb.Worksheets.Item("Recommandations")

$data1 = @(
  @("BRVM-PRINCIPAL     (**)", 0, 3, 721.97, 243.23, "🟡 Observer", "➖ Neutre"),
  @("BRVM - CONSOMMATION DE BASE     (**)", 0, 3, 719, 243.54, "🟡 Observer", "➖ Neutre"),
  @("BRVM - CONSOMMATION DISCRETIONNAIRE", 0, 3, 550.45, 188.51, "🟡 Observer", "➖ Neutre"),
  @("BRVM - INDUSTRIELS", 0, 3, 518.92, 181.2, "🟡 Observer", "➖ Neutre"),
  @("BRVM - SERVICES FINANCIERS", 0, 3, 470.45, 157.05, "🟡 Observer", "➖ Neutre"),
  @("BRVM-PRESTIGE", 0, 3, 455.44, 152.42, "🟡 Observer", "➖ Neutre"),
  @("BRVM – COMPOSITE TOTAL RETURN     (**)", 0, 3, 427.4, 143.44, "🟡 Observer", "➖ Neutre"),
  @("BRVM - ENERGIE", 0, 3, 365.21, 122.79, "🟡 Observer", "➖ Neutre"),
  @("BRVM - SERVICES PUBLICS", 0, 3, 351.25, 117.29, "🟡 Observer", "➖ Neutre"),
  @("BRVM - TELECOMMUNICATIONS", 0, 3, 295.63, 99.04000000000001, "🟡 Observer", "➖ Neutre"),
  @("EVIOSYS PACKAGING SIEM CI (SEMC)", 3, 0, 21.89, 7.3, "🟢 Achat", "✅ Renforcer"),
  @("UNILEVER CI (UNLC)", 2, 0, 14.9, 7.5, "🟡 Observer", "➖ Neutre"),
  @("UNIWAX CI (UNXC)", 2, 0, 14.61, 7.34, "🟡 Observer", "➖ Neutre"),
  @("SICABLE CI (CABC)", 2, 0, 14.17, 7.5, "🟡 Observer", "➖ Neutre"),
  @("SICOR CI (SICC)", 2, 0, 13.85, 7.47, "🟡 Observer", "➖ Neutre"),
  @("ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)", 1, 0, 7.44, 7.44, "🟡 Observer", "➖ Neutre"),
  @("SETAO CI (STAC)", 1, 0, 7.25, 7.25, "🟡 Observer", "➖ Neutre"),
  @("SITAB CI (STBC)", 1, 0, 4.97, 4.97, "🟡 Observer", "➖ Neutre"),
  @("LOTERIE NATIONALE DU BENIN (LNBB)", 0, 1, -1.02, -1.02, "🟡 Observer", "➖ Neutre"),
  @("BANK OF AFRICA NG (BOAN)", 0, 1, -1.7, -1.7, "🟡 Observer", "➖ Neutre"),
  @("CIE CI (CIEC)", 0, 1, -1.86, -1.86, "🟡 Observer", "➖ Neutre"),
  @("BANK OF AFRICA BF (BOABF)", 0, 1, -1.9, -1.9, "🟡 Observer", "➖ Neutre"),
  @("CORIS BANK INTERNATIONAL (CBIBF)", 0, 1, -2.04, -2.04, "🟡 Observer", "➖ Neutre"),
  @("ONATEL BF (ONTBF)", 0, 1, -2.21, -2.21, "🟡 Observer", "➖ Neutre"),
  @("SAPH CI (SPHC)", 0, 1, -2.4, -2.4, "🟡 Observer", "➖ Neutre"),
  @("BANK OF AFRICA ML (BOAM)", 0, 1, -2.41, -2.41, "🟡 Observer", "➖ Neutre"),
  @("ORAGROUP TOGO (ORGT)", 0, 1, -2.61, -2.61, "🟡 Observer", "➖ Neutre"),
  @("BERNABE CI (BNBC)", 0, 1, -3.25, -3.25, "🟡 Observer", "➖ Neutre"),
  @("SODE CI (SDCC)", 1, 2, -4.3, -6.89, "🟡 Observer", "👀 À surveiller"),
  @("SOGB CI (SOGC)", 0, 1, -4.34, -4.34, "🟡 Observer", "➖ Neutre"),
  @("NEI-CEDA CI (NEIC)", 0, 2, -10.59, -7.49, "🟡 Observer", "➖ Neutre")
)

for ($i = 0; $i -lt $data1.Length; $i++) {
    $row = $data1[$i]
    $r = $i + 2
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
}

# Data set shrank from 37 title rows (A2:G38) to 31 (A2:G32) - remove the leftover rows
$ws1.Range("A33:G38").EntireRow.Delete()

# --- Sheet "Top_YTD": refresh the YTD progression table (same 10 rows, new values/order) ---
$ws2 = $wb.Worksheets.Item("Top_YTD")

$data2 = @(
  @("BRVM-PRINCIPAL     (**)", 3853.04),
  @("BRVM - CONSOMMATION DE BASE     (**)", 3818.47),
  @("BRVM - CONSOMMATION DISCRETIONNAIRE", 2177.61),
  @("BRVM - INDUSTRIELS", 1932.51),
  @("BRVM - SERVICES FINANCIERS", 1593.82),
  @("BRVM-PRESTIGE", 1496.73),
  @("BRVM – COMPOSITE TOTAL RETURN     (**)", 1325.44),
  @("BRVM - ENERGIE", 990.1900000000001),
  @("BRVM - SERVICES PUBLICS", 922.97),
  @("BRVM - TELECOMMUNICATIONS", 682.64)
)

for ($i = 0; $i -lt $data2.Length; $i++) {
    $row = $data2[$i]
    $r = $i + 2
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
}

Write-Output "BRVM recommandations + Top_YTD refreshed"
